$wb = $excel.ActiveWorkbook

# The same two rows of data (for a May 1 event and a May 18 event) are
# duplicated on the "展览" sheet and the "全部类型" sheet. Both need the
# "具体时间范围" text spaced out around the dash, and the "想去人数" count
# bumped by one.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("E2").Value = "2024.05.01 09:00 - 05.01 17:00"
    $ws.Range("F2").Value = 181

    $ws.Range("E3").Value = "2024.05.18 09:00 - 05.18 17:00"
    $ws.Range("F3").Value = 124
}
